# daily auto push: 2026-01-08 09:39 UTC
# Insert a new data row for 2026/01/08 17:00 (row 593), which pushes the
# existing rows 593-634 down to 594-635 (dimension grows from D634 to D635).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 593, shifting rows 593:634 down to 594:635.
$ws.Rows.Item(593).EntireRow.Insert()

# Fill the newly inserted row with its values. Column A holds the date as
# plain text (matching every other row), so force literal text with the
# leading apostrophe and then strip the "Text" number-format style that
# gets stamped on write, so the cell matches its undecorated neighbors.
$ws.Range("A593").Value = "'2026/01/08"
$ws.Range("A593").ClearFormats()

$ws.Range("B593").Value = "木"
$ws.Range("C593").Value = 17
$ws.Range("D593").Value = 201
